$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 286, pushing the existing rows 286-407 down to 287-408.
# (All formatting/styles of row 286 carry to the new row 287, etc.)
$ws.Rows.Item(286).Insert()

# Populate the newly inserted row 286 with the new record's data.
$ws.Cells.Item(286, 1).Value  = 10
$ws.Cells.Item(286, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(286, 3).Value  = "La Araucanía"
$ws.Cells.Item(286, 4).Value  = 44875
$ws.Cells.Item(286, 5).Value  = 9
$ws.Cells.Item(286, 6).Value  = 100112009
$ws.Cells.Item(286, 7).Value  = "Acelga"
$ws.Cells.Item(286, 8).Value  = "Sin especificar"
$ws.Cells.Item(286, 9).Value  = "Primera"
$ws.Cells.Item(286, 10).Value = 215
$ws.Cells.Item(286, 11).Value = 8000
$ws.Cells.Item(286, 12).Value = 8000
$ws.Cells.Item(286, 13).Value = 8000
$ws.Cells.Item(286, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(286, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(286, 16).Value = 667
$ws.Cells.Item(286, 17).Value = 12
$ws.Cells.Item(286, 18).Value = "Hortaliza"
